$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "magapoke_2025-10-29"

# Header row: set values, then copy formatting (style) from sheet1's header
$ws2.Range("A1").Value = "rank"
$ws2.Range("B1").Value = "title"
$ws1.Range("A1:B1").Copy()
$ws2.Range("A1:B1").PasteSpecial(-4122)
$ws2.Range("A1").Select()

# Ranking data rows 2..42 (rank 1..41)
$titles = @(
    '黒月のイェルクナハト',
    'ドリーム☆ジャンボ☆ガール',
    'アイドラトリィ',
    'K-9~警視庁公安部公安第9課異能対策係~',
    '黄昏町プリズナーズ',
    '篝家の８兄弟',
    'せいぶつ部の田辺くん',
    'ハードワーカー中田',
    'ナキナギ',
    'ともだちづくり',
    'ルックスＹを選んでしまいました 〜やり込んでいるゲームに転生したはずなのに、未実装のガチャで攻略をすることになった件〜',
    '追放されなかった男　～二度目の人生は土下座から始まりました～',
    'スルガメテオ',
    '皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～',
    'お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！',
    '春くらり',
    '屋根の下のアルテミス',
    'MYS',
    '限界集落を脱村した錬金術士、都会で"最強"なのがバレまくる。～老害どもにはいい加減愛想が尽きました～',
    '生きたがりの人狼',
    'ハナバス　苔石花江のバスケ論',
    '異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～',
    '夜鐘のキト',
    'ハプスブルク家の華麗なる受難',
    'じゅーくぼっくす',
    'JK Biker',
    '平成転生',
    '鳴るさんだぁ',
    '永久のユウグレ',
    '花子狩り',
    '〈小市民〉 春期限定いちごタルト事件',
    'それがメイドのカンナです',
    '人生逆転ダンジョン',
    '卒業アルバムの彼女たち',
    '眠れる森のレガ',
    '東京デスレース',
    'ナマイキ旭ちゃんをわからせたい',
    '鉱石令嬢〜没落した悪役令嬢が炭鉱で一山当てるまでのお話〜',
    '英雄と魔女の転生ラブコメ',
    '白銀のキュイジーヌ～明治外交官の料理人～',
    'イエティ、とある日々'
)

for ($i = 0; $i -lt $titles.Length; $i++) {
    $r = $i + 2
    $ws2.Cells.Item($r, 1).Value = $i + 1
    $ws2.Cells.Item($r, 2).Value = $titles[$i]
}

Write-Host "Added sheet magapoke_2025-10-29 with $($titles.Length) ranking rows"
